# "Conclui agora o treinamento" - finish filling in the training labels
# on the "Treinamento" sheet (column B, rows 1-125) and leave the sheet
# selected/active the way the author left it (B56, scrolled near A21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

# xlCenter = -4108 ; applying it re-uses the workbook's existing "center"
# cell style (style index 3) instead of minting a new one.
$xlCenter = -4108

# B1 stays empty but picks up the centered style, matching the rest of
# column B.
$ws.Range("B1").HorizontalAlignment = $xlCenter

# B2:B5 keep their existing values but also pick up the centered style.
$ws.Range("B2").HorizontalAlignment = $xlCenter
$ws.Range("B3").HorizontalAlignment = $xlCenter
$ws.Range("B4").HorizontalAlignment = $xlCenter
$ws.Range("B5").Value = 0
$ws.Range("B5").HorizontalAlignment = $xlCenter

# B7:B125 were blank (already styled) training rows; fill in the labels.
$labels = @{
    7=0; 8=0; 9=1; 10=1; 11=0; 12=0; 13=0; 14=1; 15=1; 16=1;
    17=0; 18=1; 19=1; 20=1; 21=0; 22=0; 23=1; 24=1; 25=0; 26=1;
    27=0; 28=1; 29=0; 30=1; 31=0; 32=0; 33=0; 34=0; 35=0; 36=0;
    37=1; 38=1; 39=1; 40=0; 41=0; 42=0; 43=0; 44=0; 45=0; 46=0;
    47=0; 48=0; 49=0; 50=0; 51=0; 52=1; 53=0; 54=0; 55=1; 56=0;
    57=1; 58=1; 59=1; 60=0; 61=0; 62=0; 63=0; 64=0; 65=1; 66=0;
    67=1; 68=1; 69=1; 70=0; 71=0; 72=0; 73=0; 74=0; 75=1; 76=1;
    77=0; 78=0; 79=1; 80=1; 81=0; 82=0; 83=1; 84=0; 85=1; 86=1;
    87=0; 88=1; 89=0; 90=1; 91=1; 92=1; 93=0; 94=0; 95=1; 96=0;
    97=1; 98=1; 99=1; 100=0; 101=0; 102=0; 103=1; 104=1; 105=0; 106=1;
    107=1; 108=0; 109=1; 110=1; 111=1; 112=0; 113=0; 114=0; 115=0; 116=0;
    117=1; 118=1; 119=0; 120=1; 121=1; 122=0; 123=0; 124=1; 125=0
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 2).Value = $labels[$row]
}

# Leave "Treinamento" the active/selected sheet with B56 as the active
# cell (this also clears tabSelected on "Teste").
$ws.Activate()
$ws.Range("B56").Select()

Write-Output "done"
